$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 404-405, shifting existing rows 404-416 down to 406-418
$ws.Range("A404:R405").Insert()

# Row 404
$ws.Cells.Item(404, 1).Value = 10
$ws.Cells.Item(404, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(404, 3).Value = 'La Araucanía'
$ws.Cells.Item(404, 4).Value = 44509
$ws.Cells.Item(404, 5).Value = 9
$ws.Cells.Item(404, 6).Value = 100112003
$ws.Cells.Item(404, 7).Value = 'Ajo'
$ws.Cells.Item(404, 8).Value = 'Chino'
$ws.Cells.Item(404, 9).Value = 'Primera'
$ws.Cells.Item(404, 10).Value = 200
$ws.Cells.Item(404, 11).Value = 20000
$ws.Cells.Item(404, 12).Value = 21000
$ws.Cells.Item(404, 13).Value = 20500
$ws.Cells.Item(404, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(404, 15).Value = 'China'
$ws.Cells.Item(404, 16).Value = 2050
$ws.Cells.Item(404, 17).Value = 10
$ws.Cells.Item(404, 18).Value = 'Hortaliza'

# Row 405
$ws.Cells.Item(405, 1).Value = 10
$ws.Cells.Item(405, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(405, 3).Value = 'La Araucanía'
$ws.Cells.Item(405, 4).Value = 44509
$ws.Cells.Item(405, 5).Value = 9
$ws.Cells.Item(405, 6).Value = 100112003
$ws.Cells.Item(405, 7).Value = 'Ajo'
$ws.Cells.Item(405, 8).Value = 'Chino'
$ws.Cells.Item(405, 9).Value = 'Primera'
$ws.Cells.Item(405, 10).Value = 170
$ws.Cells.Item(405, 11).Value = 21000
$ws.Cells.Item(405, 12).Value = 22000
$ws.Cells.Item(405, 13).Value = 21471
$ws.Cells.Item(405, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(405, 15).Value = 'China'
$ws.Cells.Item(405, 16).Value = 2147
$ws.Cells.Item(405, 17).Value = 10
$ws.Cells.Item(405, 18).Value = 'Hortaliza'
